# Commit: Sat, Apr 18, 2020 10:04:57 AM
#
# Two independent changes:
#  1. The table on slide 5 gets switched from the deck's local "Table_0"
#     style to a built-in PowerPoint table style
#     ({40572E95-2EC9-403E-A9F9-04DF57A99E1C}).
#  2. The deck's theme (currently the "Integral" / Red Violet design) is
#     swapped back to the default "Office Theme" color palette. We drive
#     this the same way a user would from the Design tab: by recoloring
#     every theme color slot via the live ThemeColorScheme, which is the
#     portion of the theme actually bound to the slides/slide master.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table -------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{40572E95-2EC9-403E-A9F9-04DF57A99E1C}")
        }
    }
}

# --- 2. Recolor the theme back to the stock "Office Theme" palette ---------
# Order matches MsoThemeColorSchemeIndex / <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
